$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule R30 ("Integer min" / the "From" bound), cell C10: restore value 18 -> 1
$ws.Range("C10").Value = 1
